$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The task item that used to read "admin dashboard(Edit product)" (cell C18)
# is renamed to call out that it now also covers the orders page.
$ws.Range("C18").Value = "admin dashboard(Edit product - orders page)"

# Move / leave the active selection on the edited cell, matching the
# author's cursor position when they saved the workbook.
$ws.Range("C18").Select()
